$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ch")

# --- Update text content -------------------------------------------------
$lsquo = [char]0x2018
$ws.Range("C2").Value = "This IP only supports association with Virtual Machine or Load Balancer within the same availability zone. If you need to use multiple availability zones function, please directly switch to new BGP IP of JD Cloud"
$ws.Range("C13").Value = "$($lsquo)The resources associated."
$ws.Range("C15").Value = "$($lsquo)The Elastic Network Interface resource  associated."

# --- Formatting: red font + wrap text on the three edited cells ----------
foreach ($addr in @("C2", "C13", "C15")) {
    $cell = $ws.Range($addr)
    $cell.Font.Color = 255
    $cell.WrapText = $true
}

# --- Row heights -----------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 75
$ws.Rows.Item(15).RowHeight = 30

# --- Page setup (portrait orientation) ------------------------------------
$ws.PageSetup.Orientation = 1

# --- Selection state -----------------------------------------------------
$ws.Range("C7").Select()
